$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "26.334.46" },
    @{ Cell = "E2"; Value = "  -2.22%  " },
    @{ Cell = "D3"; Value = "1.789.11" },
    @{ Cell = "E3"; Value = "  -2.41%  " },
    @{ Cell = "D4"; Value = "'1.008" },
    @{ Cell = "E4"; Value = "  +0.22%  " },
    @{ Cell = "D5"; Value = "'1.007" },
    @{ Cell = "E5"; Value = "  +0.11%  " },
    @{ Cell = "D6"; Value = "'306.69" },
    @{ Cell = "E6"; Value = "  -1.24%  " },
    @{ Cell = "D7"; Value = "'0.4561" },
    @{ Cell = "E7"; Value = "  -1.12%  " },
    @{ Cell = "D8"; Value = "'0.3632" },
    @{ Cell = "E8"; Value = "  -0.57%  " },
    @{ Cell = "D9"; Value = "'46.27" },
    @{ Cell = "E9"; Value = "  +0.85%  " },
    @{ Cell = "D10"; Value = "'0.07077" },
    @{ Cell = "E10"; Value = "  -1.62%  " },
    @{ Cell = "D11"; Value = "'0.8742" },
    @{ Cell = "E11"; Value = "  -0.45%  " },
    @{ Cell = "D12"; Value = "'0.07799" },
    @{ Cell = "E12"; Value = "  -0.61%  " },
    @{ Cell = "D13"; Value = "'19.46" },
    @{ Cell = "E13"; Value = "  -1.01%  " },
    @{ Cell = "D14"; Value = "1.831.49" },
    @{ Cell = "E14"; Value = "  -1.25%  " },
    @{ Cell = "D15"; Value = "'5.267" },
    @{ Cell = "E15"; Value = "  -1.14%  " },
    @{ Cell = "D16"; Value = "'6.310" },
    @{ Cell = "E16"; Value = "  -1.72%  " },
    @{ Cell = "D17"; Value = "'84.78" },
    @{ Cell = "E17"; Value = "  -5.28%  " },
    @{ Cell = "D18"; Value = "'1.010" },
    @{ Cell = "E18"; Value = "  +0.36%  " },
    @{ Cell = "D19"; Value = "'0.000008514" },
    @{ Cell = "E19"; Value = "  -2.72%  " },
    @{ Cell = "D20"; Value = "'1.007" },
    @{ Cell = "E20"; Value = "  +0.12%  " },
    @{ Cell = "B21"; Value = "Avalanche" },
    @{ Cell = "C21"; Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax" },
    @{ Cell = "D21"; Value = "'14.23" },
    @{ Cell = "E21"; Value = "  -2.01%  " },
    @{ Cell = "B22"; Value = "WrappedBTC" },
    @{ Cell = "C22"; Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc" },
    @{ Cell = "D22"; Value = "26.366.98" },
    @{ Cell = "E22"; Value = "  -2.10%  " },
    @{ Cell = "D23"; Value = "'4.978" },
    @{ Cell = "E23"; Value = "  -0.34%  " },
    @{ Cell = "B24"; Value = "WrappedliquidstakedEther2.0" },
    @{ Cell = "C24"; Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth" },
    @{ Cell = "D24"; Value = "2.034.99" },
    @{ Cell = "E24"; Value = "  -3.44%  " },
    @{ Cell = "B25"; Value = "Cosmos" },
    @{ Cell = "C25"; Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom" },
    @{ Cell = "D25"; Value = "'10.50" },
    @{ Cell = "E25"; Value = "  +0.73%  " },
    @{ Cell = "D26"; Value = "'1.979" },
    @{ Cell = "E26"; Value = "  -1.16%  " },
    @{ Cell = "D27"; Value = "'152.08" },
    @{ Cell = "E27"; Value = "  +0.98%  " },
    @{ Cell = "D28"; Value = "'17.88" },
    @{ Cell = "E28"; Value = "  -1.83%  " },
    @{ Cell = "D29"; Value = "'2.031" },
    @{ Cell = "E29"; Value = "  +2.07%  " },
    @{ Cell = "D30"; Value = "'112.17" },
    @{ Cell = "E30"; Value = "  -1.73%  " },
    @{ Cell = "D31"; Value = "'4.845" },
    @{ Cell = "E31"; Value = "  -1.93%  " },
    @{ Cell = "D32"; Value = "'0.08664" },
    @{ Cell = "E32"; Value = "  -1.78%  " },
    @{ Cell = "D33"; Value = "'3.048" },
    @{ Cell = "E33"; Value = "  -1.54%  " },
    @{ Cell = "E34"; Value = "  -0.71%  " },
    @{ Cell = "D35"; Value = "'0.7189" },
    @{ Cell = "E35"; Value = "  -5.84%  " },
    @{ Cell = "D36"; Value = "'2.652" },
    @{ Cell = "E36"; Value = "  +0.54%  " },
    @{ Cell = "D37"; Value = "'1.102" },
    @{ Cell = "E37"; Value = "  -3.19%  " },
    @{ Cell = "E38"; Value = "  -0.08%  " },
    @{ Cell = "D39"; Value = "'1.077" },
    @{ Cell = "E39"; Value = "  -1.38%  " },
    @{ Cell = "D40"; Value = "'0.01943" },
    @{ Cell = "E40"; Value = "  +0.93%  " },
    @{ Cell = "D41"; Value = "'0.05096" },
    @{ Cell = "E41"; Value = "  -1.24%  " },
    @{ Cell = "D42"; Value = "'2.867" },
    @{ Cell = "E42"; Value = "  -1.79%  " },
    @{ Cell = "B43"; Value = "FraxShare" },
    @{ Cell = "C43"; Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs" },
    @{ Cell = "D43"; Value = "'6.891" },
    @{ Cell = "E43"; Value = "  -0.90%  " },
    @{ Cell = "B44"; Value = "TheSandbox" },
    @{ Cell = "C44"; Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand" },
    @{ Cell = "D44"; Value = "'0.5068" },
    @{ Cell = "E44"; Value = "  +1.50%  " },
    @{ Cell = "D45"; Value = "'0.1517" },
    @{ Cell = "E45"; Value = "  -5.08%  " },
    @{ Cell = "D46"; Value = "'8.007" },
    @{ Cell = "E46"; Value = "  -4.34%  " },
    @{ Cell = "B47"; Value = "PaxDollar" },
    @{ Cell = "C47"; Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp" },
    @{ Cell = "D47"; Value = "'1.007" },
    @{ Cell = "E47"; Value = "  +0.09%  " },
    @{ Cell = "B48"; Value = "Decentraland" },
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana" },
    @{ Cell = "D48"; Value = "'0.4649" },
    @{ Cell = "E48"; Value = "  -0.60%  " },
    @{ Cell = "D49"; Value = "'9.889" },
    @{ Cell = "E49"; Value = "  -3.58%  " },
    @{ Cell = "D50"; Value = "'99.83" },
    @{ Cell = "E50"; Value = "  -2.74%  " },
    @{ Cell = "D51"; Value = "'1.585" },
    @{ Cell = "E51"; Value = "  -1.44%  " }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
